$wb = $excel.ActiveWorkbook

# --- Worksheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 926.2353000000001
$ws.Range("J17").Value = 926.2353000000001
$ws.Range("L17").Value = 2778.7059
$ws.Range("N17").Value = -3114.7059
$ws.Range("H42").Value = 213.95
$ws.Range("I42").Value = 178.5625
$ws.Range("J42").Value = 355.5
$ws.Range("K42").Value = 535.6875
$ws.Range("L42").Value = 1066.5
$ws.Range("M42").Value = -305.6875
$ws.Range("N42").Value = -1526.5
$ws.Range("H44").Value = 20250.334
$ws.Range("I44").Value = 1001
$ws.Range("J44").Value = 29875
$ws.Range("K44").Value = 1001
$ws.Range("L44").Value = 29875
$ws.Range("M44").Value = -539
$ws.Range("N44").Value = -30799
$ws.Range("H80").Value = 1166.091
$ws.Range("I80").Value = 1094.75
$ws.Range("J80").Value = 1206.8572
$ws.Range("K80").Value = 3284.25
$ws.Range("L80").Value = 3620.5716
$ws.Range("M80").Value = -2286.25
$ws.Range("N80").Value = -5616.571599999999
$ws.Range("H83").Value = 1166.091
$ws.Range("I83").Value = 1094.75
$ws.Range("J83").Value = 1206.8572
$ws.Range("K83").Value = 9852.75
$ws.Range("L83").Value = 10861.7148
$ws.Range("M83").Value = -4860.75
$ws.Range("N83").Value = -20845.7148
$ws.Range("H88").Value = 682.93335
$ws.Range("J88").Value = 771.63635
$ws.Range("L88").Value = 771.63635
$ws.Range("N88").Value = -1583.63635
$ws.Range("H91").Value = 682.93335
$ws.Range("J91").Value = 771.63635
$ws.Range("L91").Value = 771.63635
$ws.Range("N91").Value = -3579.63635
$ws.Range("H98").Value = 3663
$ws.Range("J98").Value = 1990
$ws.Range("L98").Value = 1990
$ws.Range("N98").Value = -4986
$ws.Range("H122").Value = 3663
$ws.Range("J122").Value = 1990
$ws.Range("L122").Value = 5970
$ws.Range("N122").Value = -10870
$ws.Range("H137").Value = 5005912
$ws.Range("I137").Value = 10872561
$ws.Range("J137").Value = 8395.814
$ws.Range("K137").Value = 32617683
$ws.Range("L137").Value = 25187.442
$ws.Range("M137").Value = -32615133
$ws.Range("N137").Value = -30287.442
$ws.Range("H138").Value = 6680
$ws.Range("J138").Value = 7475
$ws.Range("L138").Value = 22425
$ws.Range("N138").Value = -32705

# --- Worksheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 467.42856
$ws.Range("I4").Value = 437
$ws.Range("J4").Value = 579
$ws.Range("K4").Value = 437
$ws.Range("L4").Value = 579
$ws.Range("M4").Value = -321
$ws.Range("N4").Value = -811
$ws.Range("H32").Value = 3395.7925
$ws.Range("I32").Value = 3097.3953
$ws.Range("J32").Value = 4678.9
$ws.Range("K32").Value = 3097.3953
$ws.Range("L32").Value = 4678.9
$ws.Range("M32").Value = -2810.3953
$ws.Range("N32").Value = -5252.9
$ws.Range("H74").Value = 2386.0967
$ws.Range("I74").Value = 796.1579
$ws.Range("J74").Value = 4903.5
$ws.Range("K74").Value = 796.1579
$ws.Range("L74").Value = 4903.5
$ws.Range("M74").Value = 77.84209999999996
$ws.Range("N74").Value = -6651.5
$ws.Range("H77").Value = 2386.0967
$ws.Range("I77").Value = 796.1579
$ws.Range("J77").Value = 4903.5
$ws.Range("K77").Value = 3980.7895
$ws.Range("L77").Value = 24517.5
$ws.Range("M77").Value = 387.2104999999997
$ws.Range("N77").Value = -33253.5
$ws.Range("H112").Value = 37999
$ws.Range("J112").Value = 37999
$ws.Range("L112").Value = 37999
$ws.Range("N112").Value = -40953
$ws.Range("H132").Value = 5260.0684
$ws.Range("I132").Value = 2627.2903
$ws.Range("K132").Value = 7881.8709
$ws.Range("M132").Value = -5351.8709

# --- Worksheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H56").Value = 51994.668
$ws.Range("I56").Value = 51994.668
$ws.Range("K56").Value = 51994.668
$ws.Range("M56").Value = -51255.668
$ws.Range("H86").Value = 3322
$ws.Range("I86").Value = 2755.25
$ws.Range("K86").Value = 2755.25
$ws.Range("M86").Value = -1632.25
$ws.Range("H89").Value = 3322
$ws.Range("I89").Value = 2755.25
$ws.Range("K89").Value = 13776.25
$ws.Range("M89").Value = -8160.25

# --- Worksheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 17437.5
$ws.Range("J41").Value = 60000
$ws.Range("L41").Value = 60000
$ws.Range("N41").Value = -60856
$ws.Range("H62").Value = 19444.6
$ws.Range("J62").Value = 26666.5
$ws.Range("L62").Value = 26666.5
$ws.Range("N62").Value = -27914.5
$ws.Range("H65").Value = 19444.6
$ws.Range("J65").Value = 26666.5
$ws.Range("L65").Value = 133332.5
$ws.Range("N65").Value = -139572.5
$ws.Range("H74").Value = 46666.25
$ws.Range("J74").Value = 47999.5
$ws.Range("L74").Value = 47999.5
$ws.Range("N74").Value = -49747.5
$ws.Range("H77").Value = 46666.25
$ws.Range("J77").Value = 47999.5
$ws.Range("L77").Value = 143998.5
$ws.Range("N77").Value = -152734.5
$ws.Range("H107").Value = 2327.3333
$ws.Range("I107").Value = 1599.9166
$ws.Range("K107").Value = 1599.9166
$ws.Range("M107").Value = 320.0834
$ws.Range("H132").Value = 46819.766
$ws.Range("I132").Value = 1999.9
$ws.Range("J132").Value = 136459.5
$ws.Range("K132").Value = 5999.700000000001
$ws.Range("L132").Value = 409378.5
$ws.Range("M132").Value = -3469.700000000001
$ws.Range("N132").Value = -414438.5

# --- Worksheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H25").Value = 512.75
$ws.Range("I25").Value = 27
$ws.Range("K25").Value = 81
$ws.Range("M25").Value = 88
$ws.Range("H30").Value = 512.75
$ws.Range("I30").Value = 27
$ws.Range("K30").Value = 81
$ws.Range("M30").Value = 21
$ws.Range("H68").Value = 98005.234
$ws.Range("J68").Value = 3401.2
$ws.Range("L68").Value = 10203.6
$ws.Range("N68").Value = -11825.6
$ws.Range("H71").Value = 98005.234
$ws.Range("J71").Value = 3401.2
$ws.Range("L71").Value = 30610.8
$ws.Range("N71").Value = -38722.8
$ws.Range("H113").Value = 3472.4546
$ws.Range("J113").Value = 3931
$ws.Range("L113").Value = 11793
$ws.Range("N113").Value = -16133

# --- Worksheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H54").Value = 7566.6665
$ws.Range("J54").Value = 7250
$ws.Range("L54").Value = 7250
$ws.Range("N54").Value = -8030
$ws.Range("H70").Value = 7685.6665
$ws.Range("I70").Value = 7524.5835
$ws.Range("K70").Value = 7524.5835
$ws.Range("M70").Value = -7254.5835
$ws.Range("H73").Value = 7685.6665
$ws.Range("I73").Value = 7524.5835
$ws.Range("K73").Value = 7524.5835
$ws.Range("M73").Value = -6588.5835
$ws.Range("H122").Value = 3694.2068
$ws.Range("I122").Value = 3985.28
$ws.Range("J122").Value = 1875
$ws.Range("K122").Value = 11955.84
$ws.Range("L122").Value = 5625
$ws.Range("M122").Value = -9505.84
$ws.Range("N122").Value = -10525
$ws.Range("H123").Value = 36929.223
$ws.Range("J123").Value = 36929.223
$ws.Range("L123").Value = 36929.223
$ws.Range("N123").Value = -41829.223

# --- Worksheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 6887.5884
$ws.Range("I46").Value = 1598
$ws.Range("K46").Value = 1598
$ws.Range("M46").Value = -1410
$ws.Range("H55").Value = 1079.3889
$ws.Range("I55").Value = 875.3333
$ws.Range("K55").Value = 875.3333
$ws.Range("M55").Value = -702.3333
$ws.Range("H61").Value = 4383.1665
$ws.Range("J61").Value = 5899.5
$ws.Range("L61").Value = 5899.5
$ws.Range("N61").Value = -6303.5
$ws.Range("H113").Value = 4383.1665
$ws.Range("J113").Value = 5899.5
$ws.Range("L113").Value = 5899.5
$ws.Range("N113").Value = -10239.5
$ws.Range("H132").Value = 4021.5227
$ws.Range("I132").Value = 3672.8
$ws.Range("J132").Value = 4312.125
$ws.Range("K132").Value = 11018.4
$ws.Range("L132").Value = 12936.375
$ws.Range("M132").Value = -8488.400000000001
$ws.Range("N132").Value = -17996.375
$ws.Range("H133").Value = 48938.668
$ws.Range("J133").Value = 48938.668
$ws.Range("L133").Value = 48938.668
$ws.Range("N133").Value = -53998.668
$ws.Range("H136").Value = 9449.6875
$ws.Range("I136").Value = 7439
$ws.Range("K136").Value = 22317
$ws.Range("M136").Value = -19767
